$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 7: Bowling Green --------------------------------------------
$ws.Range("A7").Value = "Bowling Green"
$ws.Range("B7").Value = "Ohio"
$ws.Range("C7").Value = 185
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = "bgsu.edu"

# Hyperlink for the new row's Website cell (added after the text so the
# literal cell text is preserved instead of being overwritten by the
# hyperlink's display text).
$ws.Hyperlinks.Add($ws.Range("E7"), "https://www.cam.ac.uk/", [Type]::Missing, [Type]::Missing, "https://www.cam.ac.uk/")
$ws.Range("E7").Value = "bgsu.edu"
# Match the rest of the sheet: hyperlinked cells keep the plain left/top
# style (no blue/underline "Hyperlink" style) - restore the same formatting
# used by the rest of the row/table by copying it from a sibling cell.
$ws.Range("A7").Copy() | Out-Null
$ws.Range("E7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- New column F: "Has Hospital" ---------------------------------------
$ws.Range("F2").Value = "Yes"
$ws.Range("F3").Value = "No"
$ws.Range("F4").Value = "yes"
$ws.Range("F5").Value = "No"
$ws.Range("F6").Value = "yes"
$ws.Range("F7").Value = "no"

$ws.Range("F1").Value = "Has Hospital"
# Reset to the workbook default style first so the font change below does
# not also carry along the column's inherited left/top alignment.
$ws.Range("F1").Style = "Normal"
$f1Font = $ws.Range("F1").Font
$f1Font.Size = 10
$f1Font.Name = "Arial"

# --- Selection -------------------------------------------------------------
$ws.Range("F1").Select() | Out-Null
